$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update K3: new comment text, change fill from yellow to light gray (same style as rows 4/5)
$ws.Range("K3").Value = "pretty darn good…I think so"
$ws.Range("A4").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# restore selection to K3 (matches recorded selection change)
$ws.Range("K3").Select()
